$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 7).Value = 21.267222
$ws.Cells.Item(2, 8).Value = 63.801666
$ws.Cells.Item(2, 9).Value = 0.06271644651145813
$ws.Cells.Item(2, 10).Value = 0.06271644651145813
$ws.Cells.Item(2, 13).Value = 0.467036
$ws.Cells.Item(2, 14).Value = 1.401108
$ws.Cells.Item(2, 15).Value = 0.001972893265924874
$ws.Cells.Item(2, 16).Value = 0.001972893265924874
$ws.Cells.Item(2, 17).Value = 9.932558293992001
$ws.Cells.Item(2, 18).Value = 89.39302464592799
$ws.Cells.Item(2, 19).Value = 0.0001237328549851933
$ws.Cells.Item(2, 20).Value = 0.0001237328549851933
$ws.Cells.Item(3, 7).Value = 21.267222
$ws.Cells.Item(3, 8).Value = 63.801666
$ws.Cells.Item(3, 9).Value = 0.06271644651145813
$ws.Cells.Item(3, 10).Value = 0.06271644651145813
$ws.Cells.Item(3, 15).Value = 0.001374344438283074
$ws.Cells.Item(3, 16).Value = 0.001374344438283074
$ws.Cells.Item(3, 17).Value = 6.919155985293999
$ws.Cells.Item(3, 18).Value = 62.272403867646
$ws.Cells.Item(3, 19).Value = 0.00008619399945190038
$ws.Cells.Item(3, 20).Value = 0.00008619399945190038
$ws.Cells.Item(4, 7).Value = 21.267222
$ws.Cells.Item(4, 8).Value = 63.801666
$ws.Cells.Item(4, 9).Value = 0.06271644651145813
$ws.Cells.Item(4, 10).Value = 0.06271644651145813
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.02089333333333333
$ws.Cells.Item(4, 14).Value = 0.06268
$ws.Cells.Item(4, 15).Value = 0.00008825939892440207
$ws.Cells.Item(4, 16).Value = 0.00008825939892440207
$ws.Cells.Item(4, 17).Value = 0.44434315832
$ws.Cells.Item(4, 18).Value = 3.99908842488
$ws.Cells.Item(4, 19).Value = 0.000005535315871775708
$ws.Cells.Item(4, 20).Value = 0.000005535315871775708
$ws.Cells.Item(5, 7).Value = 21.267222
$ws.Cells.Item(5, 8).Value = 63.801666
$ws.Cells.Item(5, 9).Value = 0.06271644651145813
$ws.Cells.Item(5, 10).Value = 0.06271644651145813
$ws.Cells.Item(5, 13).Value = 235.9131673333333
$ws.Cells.Item(5, 14).Value = 707.739502
$ws.Cells.Item(5, 15).Value = 0.9965645028968676
$ws.Cells.Item(5, 16).Value = 0.9965645028968676
$ws.Cells.Item(5, 17).Value = 5017.217702401148
$ws.Cells.Item(5, 18).Value = 45154.95932161033
$ws.Cells.Item(5, 19).Value = 0.06250098434114926
$ws.Cells.Item(5, 20).Value = 0.06250098434114926
$ws.Cells.Item(6, 9).Value = 0.4054090708715844
$ws.Cells.Item(6, 10).Value = 0.4054090708715843
$ws.Cells.Item(6, 13).Value = 0.467036
$ws.Cells.Item(6, 14).Value = 1.401108
$ws.Cells.Item(6, 15).Value = 0.001972893265924874
$ws.Cells.Item(6, 16).Value = 0.001972893265924874
$ws.Cells.Item(6, 17).Value = 64.20563430055734
$ws.Cells.Item(6, 18).Value = 577.850708705016
$ws.Cells.Item(6, 19).Value = 0.000799828825867409
$ws.Cells.Item(6, 20).Value = 0.0007998288258674089
$ws.Cells.Item(7, 9).Value = 0.4054090708715844
$ws.Cells.Item(7, 10).Value = 0.4054090708715843
$ws.Cells.Item(7, 15).Value = 0.001374344438283074
$ws.Cells.Item(7, 16).Value = 0.001374344438283074
$ws.Cells.Item(7, 19).Value = 0.0005571717017818705
$ws.Cells.Item(7, 20).Value = 0.0005571717017818704
$ws.Cells.Item(8, 9).Value = 0.4054090708715844
$ws.Cells.Item(8, 10).Value = 0.4054090708715843
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.02089333333333333
$ws.Cells.Item(8, 14).Value = 0.06268
$ws.Cells.Item(8, 15).Value = 0.00008825939892440207
$ws.Cells.Item(8, 16).Value = 0.00008825939892440207
$ws.Cells.Item(8, 17).Value = 2.872304745928889
$ws.Cells.Item(8, 18).Value = 25.85074271336
$ws.Cells.Item(8, 19).Value = 0.00003578116091362636
$ws.Cells.Item(8, 20).Value = 0.00003578116091362635
$ws.Cells.Item(9, 9).Value = 0.4054090708715844
$ws.Cells.Item(9, 10).Value = 0.4054090708715843
$ws.Cells.Item(9, 13).Value = 235.9131673333333
$ws.Cells.Item(9, 14).Value = 707.739502
$ws.Cells.Item(9, 15).Value = 0.9965645028968676
$ws.Cells.Item(9, 16).Value = 0.9965645028968676
$ws.Cells.Item(9, 17).Value = 32432.09206247525
$ws.Cells.Item(9, 18).Value = 291888.8285622772
$ws.Cells.Item(9, 19).Value = 0.4040162891830215
$ws.Cells.Item(9, 20).Value = 0.4040162891830214
$ws.Cells.Item(10, 7).Value = 121.820091
$ws.Cells.Item(10, 8).Value = 365.460273
$ws.Cells.Item(10, 9).Value = 0.3592440621169263
$ws.Cells.Item(10, 10).Value = 0.3592440621169263
$ws.Cells.Item(10, 13).Value = 0.467036
$ws.Cells.Item(10, 14).Value = 1.401108
$ws.Cells.Item(10, 15).Value = 0.001972893265924874
$ws.Cells.Item(10, 16).Value = 0.001972893265924874
$ws.Cells.Item(10, 17).Value = 56.894368020276
$ws.Cells.Item(10, 18).Value = 512.0493121824841
$ws.Cells.Item(10, 19).Value = 0.0007087501909739813
$ws.Cells.Item(10, 20).Value = 0.0007087501909739813
$ws.Cells.Item(11, 7).Value = 121.820091
$ws.Cells.Item(11, 8).Value = 365.460273
$ws.Cells.Item(11, 9).Value = 0.3592440621169263
$ws.Cells.Item(11, 10).Value = 0.3592440621169263
$ws.Cells.Item(11, 15).Value = 0.001374344438283074
$ws.Cells.Item(11, 16).Value = 0.001374344438283074
$ws.Cells.Item(11, 17).Value = 39.633395079607
$ws.Cells.Item(11, 18).Value = 356.700555716463
$ws.Cells.Item(11, 19).Value = 0.0004937250787566169
$ws.Cells.Item(11, 20).Value = 0.0004937250787566169
$ws.Cells.Item(12, 7).Value = 121.820091
$ws.Cells.Item(12, 8).Value = 365.460273
$ws.Cells.Item(12, 9).Value = 0.3592440621169263
$ws.Cells.Item(12, 10).Value = 0.3592440621169263
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.02089333333333333
$ws.Cells.Item(12, 14).Value = 0.06268
$ws.Cells.Item(12, 15).Value = 0.00008825939892440207
$ws.Cells.Item(12, 16).Value = 0.00008825939892440207
$ws.Cells.Item(12, 17).Value = 2.54522776796
$ws.Cells.Item(12, 18).Value = 22.90704991164
$ws.Cells.Item(12, 19).Value = 0.00003170666498960048
$ws.Cells.Item(12, 20).Value = 0.00003170666498960048
$ws.Cells.Item(13, 7).Value = 121.820091
$ws.Cells.Item(13, 8).Value = 365.460273
$ws.Cells.Item(13, 9).Value = 0.3592440621169263
$ws.Cells.Item(13, 10).Value = 0.3592440621169263
$ws.Cells.Item(13, 13).Value = 235.9131673333333
$ws.Cells.Item(13, 14).Value = 707.739502
$ws.Cells.Item(13, 15).Value = 0.9965645028968676
$ws.Cells.Item(13, 16).Value = 0.9965645028968676
$ws.Cells.Item(13, 17).Value = 28738.9635126449
$ws.Cells.Item(13, 18).Value = 258650.6716138041
$ws.Cells.Item(13, 19).Value = 0.3580098801822061
$ws.Cells.Item(13, 20).Value = 0.3580098801822061
$ws.Cells.Item(14, 7).Value = 58.539182
$ws.Cells.Item(14, 8).Value = 175.617546
$ws.Cells.Item(14, 9).Value = 0.1726304205000311
$ws.Cells.Item(14, 10).Value = 0.1726304205000311
$ws.Cells.Item(14, 13).Value = 0.467036
$ws.Cells.Item(14, 14).Value = 1.401108
$ws.Cells.Item(14, 15).Value = 0.001972893265924874
$ws.Cells.Item(14, 16).Value = 0.001972893265924874
$ws.Cells.Item(14, 17).Value = 27.339905404552
$ws.Cells.Item(14, 18).Value = 246.059148640968
$ws.Cells.Item(14, 19).Value = 0.0003405813940982908
$ws.Cells.Item(14, 20).Value = 0.0003405813940982907
$ws.Cells.Item(15, 7).Value = 58.539182
$ws.Cells.Item(15, 8).Value = 175.617546
$ws.Cells.Item(15, 9).Value = 0.1726304205000311
$ws.Cells.Item(15, 10).Value = 0.1726304205000311
$ws.Cells.Item(15, 15).Value = 0.001374344438283074
$ws.Cells.Item(15, 16).Value = 0.001374344438283074
$ws.Cells.Item(15, 17).Value = 19.04535211554733
$ws.Cells.Item(15, 18).Value = 171.408169039926
$ws.Cells.Item(15, 19).Value = 0.0002372536582926861
$ws.Cells.Item(15, 20).Value = 0.0002372536582926861
$ws.Cells.Item(16, 7).Value = 58.539182
$ws.Cells.Item(16, 8).Value = 175.617546
$ws.Cells.Item(16, 9).Value = 0.1726304205000311
$ws.Cells.Item(16, 10).Value = 0.1726304205000311
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.02089333333333333
$ws.Cells.Item(16, 14).Value = 0.06268
$ws.Cells.Item(16, 15).Value = 0.00008825939892440207
$ws.Cells.Item(16, 16).Value = 0.00008825939892440207
$ws.Cells.Item(16, 17).Value = 1.223078642586667
$ws.Cells.Item(16, 18).Value = 11.00770778328
$ws.Cells.Item(16, 19).Value = 0.00001523625714939952
$ws.Cells.Item(16, 20).Value = 0.00001523625714939952
$ws.Cells.Item(17, 7).Value = 58.539182
$ws.Cells.Item(17, 8).Value = 175.617546
$ws.Cells.Item(17, 9).Value = 0.1726304205000311
$ws.Cells.Item(17, 10).Value = 0.1726304205000311
$ws.Cells.Item(17, 13).Value = 235.9131673333333
$ws.Cells.Item(17, 14).Value = 707.739502
$ws.Cells.Item(17, 15).Value = 0.9965645028968676
$ws.Cells.Item(17, 16).Value = 0.9965645028968676
$ws.Cells.Item(17, 17).Value = 13810.16383872246
$ws.Cells.Item(17, 18).Value = 124291.4745485021
$ws.Cells.Item(17, 19).Value = 0.1720373491904907
$ws.Cells.Item(17, 20).Value = 0.1720373491904907
